$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 360797.12
$ws.Range("J17").Value = 375804.38
$ws.Range("L17").Value = 1127413.14
$ws.Range("N17").Value = -1127749.14

$ws.Range("H18").Value = 1787.5
$ws.Range("I18").Value = 2116.6667
$ws.Range("J18").Value = 800
$ws.Range("K18").Value = 2116.6667
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = -1832.6667
$ws.Range("N18").Value = -1368

$ws.Range("H19").Value = 468.3
$ws.Range("I19").Value = 480.2
$ws.Range("J19").Value = 456.4
$ws.Range("K19").Value = 480.2
$ws.Range("L19").Value = 456.4
$ws.Range("M19").Value = -305.2
$ws.Range("N19").Value = -806.4

$ws.Range("H74").Value = 3397.95
$ws.Range("I74").Value = 3326.6
$ws.Range("J74").Value = 3612
$ws.Range("K74").Value = 3326.6
$ws.Range("L74").Value = 3612
$ws.Range("M74").Value = -2390.6
$ws.Range("N74").Value = -5484

$ws.Range("H77").Value = 3397.95
$ws.Range("I77").Value = 3326.6
$ws.Range("J77").Value = 3612
$ws.Range("K77").Value = 16633
$ws.Range("L77").Value = 18060
$ws.Range("M77").Value = -11953
$ws.Range("N77").Value = -27420

$ws.Range("H101").Value = 7037
$ws.Range("I101").Value = 1000
$ws.Range("J101").Value = 8546.25
$ws.Range("K101").Value = 3000
$ws.Range("L101").Value = 25638.75
$ws.Range("M101").Value = -1378
$ws.Range("N101").Value = -28882.75

$ws.Range("H137").Value = 2273.125
$ws.Range("I137").Value = 2132.8572
$ws.Range("J137").Value = 2731.3333
$ws.Range("K137").Value = 6398.571599999999
$ws.Range("L137").Value = 8193.999899999999
$ws.Range("M137").Value = -3848.571599999999
$ws.Range("N137").Value = -13293.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4610.0293
$ws.Range("I32").Value = 3519.322
$ws.Range("J32").Value = 11760.223
$ws.Range("K32").Value = 3519.322
$ws.Range("L32").Value = 11760.223
$ws.Range("M32").Value = -3232.322
$ws.Range("N32").Value = -12334.223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 38448
$ws.Range("J93").Value = 38448
$ws.Range("L93").Value = 38448
$ws.Range("N93").Value = -42192

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 829.375
$ws.Range("I22").Value = 951.53845
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 951.53845
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = -601.53845
$ws.Range("N22").Value = -1000

$ws.Range("H58").Value = 1020.3571
$ws.Range("I58").Value = 1018.3333
$ws.Range("J58").Value = 1028.6364
$ws.Range("K58").Value = 1018.3333
$ws.Range("L58").Value = 1028.6364
$ws.Range("M58").Value = -815.3333
$ws.Range("N58").Value = -1434.6364

$ws.Range("H99").Value = 1872.2354
$ws.Range("I99").Value = 1625
$ws.Range("J99").Value = 2465.6
$ws.Range("K99").Value = 1625
$ws.Range("L99").Value = 2465.6
$ws.Range("M99").Value = -127
$ws.Range("N99").Value = -5461.6

$ws.Range("H126").Value = 1872.2354
$ws.Range("I126").Value = 1625
$ws.Range("J126").Value = 2465.6
$ws.Range("K126").Value = 4875
$ws.Range("L126").Value = 7396.799999999999
$ws.Range("M126").Value = -2405
$ws.Range("N126").Value = -12336.8

$ws.Range("H136").Value = 1020.3571
$ws.Range("I136").Value = 1018.3333
$ws.Range("J136").Value = 1028.6364
$ws.Range("K136").Value = 3054.9999
$ws.Range("L136").Value = 3085.9092
$ws.Range("M136").Value = -504.9998999999998
$ws.Range("N136").Value = -8185.9092

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 801.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 801.5
$ws.Range("K92").Value = 0
$ws.Range("M92").Value = 2404.5
$ws.Range("N92").Value = -4900.5
$ws.Range("L92").ClearContents()

$ws.Range("H107").Value = 541.8889
$ws.Range("J107").Value = 720.8
$ws.Range("L107").Value = 2162.4
$ws.Range("N107").Value = -6002.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1777.6923
$ws.Range("I7").Value = 1460
$ws.Range("J7").Value = 2836.6667
$ws.Range("K7").Value = 1460
$ws.Range("L7").Value = 2836.6667
$ws.Range("M7").Value = -1348
$ws.Range("N7").Value = -3060.6667

$ws.Range("H40").Value = 2936.9092
$ws.Range("J40").Value = 2684.0908
$ws.Range("L40").Value = 2684.0908
$ws.Range("N40").Value = -2956.0908

$ws.Range("H55").Value = 644.1111
$ws.Range("I55").Value = 579.8
$ws.Range("K55").Value = 579.8
$ws.Range("M55").Value = -406.8

$ws.Range("H61").Value = 1811.3334
$ws.Range("I61").Value = 1881.0769
$ws.Range("J61").Value = 1630
$ws.Range("K61").Value = 1881.0769
$ws.Range("L61").Value = 1630
$ws.Range("M61").Value = -1679.0769
$ws.Range("N61").Value = -2034

$ws.Range("H113").Value = 1811.3334
$ws.Range("I113").Value = 1881.0769
$ws.Range("J113").Value = 1630
$ws.Range("K113").Value = 1881.0769
$ws.Range("L113").Value = 1630
$ws.Range("M113").Value = 288.9231
$ws.Range("N113").Value = -5970

$ws.Range("H122").Value = 4264.077
$ws.Range("I122").Value = 6977
$ws.Range("J122").Value = 3058.3333
$ws.Range("K122").Value = 20931
$ws.Range("L122").Value = 9174.999899999999
$ws.Range("M122").Value = -18481
$ws.Range("N122").Value = -14074.9999

$ws.Range("H126").Value = 1777.6923
$ws.Range("I126").Value = 1460
$ws.Range("J126").Value = 2836.6667
$ws.Range("K126").Value = 4380
$ws.Range("L126").Value = 8510.000100000001
$ws.Range("M126").Value = -1910
$ws.Range("N126").Value = -13450.0001

$ws.Range("H136").Value = 4412.569
$ws.Range("I136").Value = 2135.853
$ws.Range("J136").Value = 7637.9165
$ws.Range("K136").Value = 6407.559
$ws.Range("L136").Value = 22913.7495
$ws.Range("M136").Value = -3857.559
$ws.Range("N136").Value = -28013.7495

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2977.4783
$ws.Range("I81").Value = 578
$ws.Range("J81").Value = 14375
$ws.Range("K81").Value = 1156
$ws.Range("L81").Value = 28750
$ws.Range("M81").Value = -95
$ws.Range("N81").Value = -30872

$ws.Range("H84").Value = 2977.4783
$ws.Range("I84").Value = 578
$ws.Range("J84").Value = 14375
$ws.Range("K84").Value = 5780
$ws.Range("L84").Value = 143750
$ws.Range("M84").Value = -476
$ws.Range("N84").Value = -154358

$ws.Range("H132").Value = 1258.686
$ws.Range("I132").Value = 1039.4706
$ws.Range("J132").Value = 2086.8333
$ws.Range("K132").Value = 3118.4118
$ws.Range("L132").Value = 6260.499899999999
$ws.Range("M132").Value = -588.4118000000003
$ws.Range("N132").Value = -11320.4999
